$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary figures ---
# Valor Mora (E11): 400228 -> 118036
$ws.Range("E11").Value = 118036

# Cant. Trabajadores (C13): 8 -> 2 ; Cant. Periodos (F13): 9 -> 4
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 4

# --- Replace the first data row (row 16) with the new worker/record ---
$ws.Range("C16").Value = "87246999"
$ws.Range("D16").Value = "RUBIO RAIMUNDO BOLAÑOS LOPEZ"
$ws.Range("E16").Value = "1703"
$ws.Range("F16").Value = 29509
$ws.Range("G16").Value = 737717

# --- Update the periods for the GUSTAVO ENRIQUE CALVO MORE rows (17 & 18) ---
$ws.Range("E17").Value = "1706"
$ws.Range("E18").Value = "1707"

# --- Copy the formatting (last-row borders) from the old closing data row (27)
#     onto row 19, which becomes the new last row of the table ---
$ws.Range("B27:J27").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# Give row 19 the 4th (last) worker record
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143400321"
$ws.Range("D19").Value = "GUSTAVO ENRIQUE CALVO MORE"
$ws.Range("E19").Value = "1708"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 877803
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = ""
$ws.Range("J19").Value = ""

# --- Remove the now-obsolete data rows (old rows 20-27); this also shifts the
#     signature block (old rows 32/33) up to rows 24/25 ---
$ws.Rows("20:27").Delete()

# --- Column D is now narrower because the longest name shrank ---
$ws.Columns("D:D").AutoFit()
